# Insert a new data row at row 445 (pushing the existing rows 445-501 down
# to 446-502) and populate it with the new "Ajo" price record for
# 2023-08-16 (serial date 45154).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(445).Insert()

$ws.Cells.Item(445, 1).Value = 7
$ws.Cells.Item(445, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(445, 3).Value = "Ñuble"
$ws.Cells.Item(445, 4).Value = 45154
$ws.Cells.Item(445, 5).Value = 16
$ws.Cells.Item(445, 6).Value = 100112003
$ws.Cells.Item(445, 7).Value = "Ajo"
$ws.Cells.Item(445, 8).Value = "Chino"
$ws.Cells.Item(445, 9).Value = "Primera"
$ws.Cells.Item(445, 10).Value = 30
$ws.Cells.Item(445, 11).Value = 21000
$ws.Cells.Item(445, 12).Value = 21000
$ws.Cells.Item(445, 13).Value = 21000
$ws.Cells.Item(445, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(445, 15).Value = "China"
$ws.Cells.Item(445, 16).Value = 2100
$ws.Cells.Item(445, 17).Value = 10
$ws.Cells.Item(445, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(445, 4).NumberFormat = $ws.Cells.Item(446, 4).NumberFormat
